# Roving Diver Code table update: change PERMANOVA factor "County" -> "Location"
# (site-level analysis), and resize the second table column / affected rows
# to match the new, shorter label widths.

$d = $word.ActiveDocument

# Replace the "County" factor label wherever it appears as a whole word,
# including inside "County:Date" (':' is a word boundary), turning it into
# "Location" / "Location:Date".
$d.Content.Find.Execute("County", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Location", 2)

# The results table is the first (only) table in the document.
$t = $d.Tables.Item(1)

# Narrower "Location" text allows the second column to shrink slightly.
$t.Columns.Item(2).Width = 88.4

# The two rows whose labels changed become a touch shorter as well.
$t.Rows.Item(2).Height = 20.7
$t.Rows.Item(4).Height = 20.55
